$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column F (dSF) values for the affected rows
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = -5
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = -6
$ws.Range("F17").Value = 0
